# Sample Project / Main.xlsx - "Rules" sheet
# Cell B11 (the 4th rule row, "R40") is retyped as the text "1".
# A leading apostrophe forces Excel to store the numeric-looking entry
# as text (a shared string) instead of silently converting it to the
# number 1, matching the sharedStrings.xml addition of a new <si><t>1</t>
# entry and B11's <c t="s"> pointing at it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").Value = "'1"
